$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text  = "61÷5=12, 1"
$tbl.Cell(1,2).Range.Text  = "81÷5=16, 1"
$tbl.Cell(1,3).Range.Text  = "12÷5=2, 2"
$tbl.Cell(1,4).Range.Text  = "66÷6=11, 0"
$tbl.Cell(1,5).Range.Text  = "86÷9=9, 5"

$tbl.Cell(5,1).Range.Text  = "33÷6=5, 3"
$tbl.Cell(5,2).Range.Text  = "16÷6=2, 4"
$tbl.Cell(5,3).Range.Text  = "42÷4=10, 2"
$tbl.Cell(5,4).Range.Text  = "47÷7=6, 5"
$tbl.Cell(5,5).Range.Text  = "12÷3=4, 0"

$tbl.Cell(9,1).Range.Text  = "65÷5=13, 0"
$tbl.Cell(9,2).Range.Text  = "63÷5=12, 3"
$tbl.Cell(9,3).Range.Text  = "56÷6=9, 2"
$tbl.Cell(9,4).Range.Text  = "97÷6=16, 1"
$tbl.Cell(9,5).Range.Text  = "31÷7=4, 3"

$tbl.Cell(13,1).Range.Text = "57÷8=7, 1"
$tbl.Cell(13,2).Range.Text = "98÷9=10, 8"
$tbl.Cell(13,3).Range.Text = "97÷6=16, 1"
$tbl.Cell(13,4).Range.Text = "55÷2=27, 1"
$tbl.Cell(13,5).Range.Text = "56÷5=11, 1"

$tbl.Cell(17,1).Range.Text = "77÷4=19, 1"
$tbl.Cell(17,2).Range.Text = "94÷4=23, 2"
$tbl.Cell(17,3).Range.Text = "42÷7=6, 0"
$tbl.Cell(17,4).Range.Text = "96÷5=19, 1"
$tbl.Cell(17,5).Range.Text = "69÷2=34, 1"
